$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.616.66'
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").Value = '1.588.00'
$ws.Range("E3").Value = '  -2.59%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -2.64%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("E12").Value = '  -2.70%  '
$ws.Range("D13").Value = '1.592.07'
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '26.616.75'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '207.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("E22").Value = '  -3.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.08%  '
$ws.Range("E24").Value = '  -2.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("E32").Value = '  -3.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.664'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +22.92%  '
$ws.Range("D34").Value = '1.327.34'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.97%  '
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.40'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.786'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("E43").Value = '  -3.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '1.722.58'
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.831'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0976'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.88%  '
